$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 183; this shifts the existing rows 183:222
# down to 184:223, preserving their original content and formatting.
$ws.Rows(183).Insert()

# Populate the newly inserted row 183 with the new weekly price record.
$ws.Cells.Item(183, 1).Value = 4
$ws.Cells.Item(183, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(183, 3).Value = "Los Lagos"
$ws.Cells.Item(183, 4).Value = (Get-Date -Year 2021 -Month 12 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(183, 5).Value = 10
$ws.Cells.Item(183, 6).Value = 100112045
$ws.Cells.Item(183, 7).Value = "Zapallo"
$ws.Cells.Item(183, 8).Value = "Camote"
$ws.Cells.Item(183, 9).Value = "1a nueva(o)"
$ws.Cells.Item(183, 10).Value = 100
$ws.Cells.Item(183, 11).Value = 600
$ws.Cells.Item(183, 12).Value = 600
$ws.Cells.Item(183, 13).Value = 600
$ws.Cells.Item(183, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(183, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(183, 16).Value = 600
$ws.Cells.Item(183, 17).Value = 1
$ws.Cells.Item(183, 18).Value = "Hortaliza"
